$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.244.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.886.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.96%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.82"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.688"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.85"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +8.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.350"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.70"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.11%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0970"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.09"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.164.02"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.723"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.903.16"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.227.50"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.00"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.84%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "245.32"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.81"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.96"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.31%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.94%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.35%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -10.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.01"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.48"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.29"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.22%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128.46"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.76"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +8.31%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0579"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.22"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.67%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.851"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.23%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.85%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -23.16%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.59"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0669"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.28%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.20%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.292.71"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.30%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0808"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +8.44%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.40"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.53%  "
